# Update cryptocurrency Price (D) and Volume(1h) (E) cells with refreshed market data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $value) {
    # Force the cell to remain plain text (matches the source data, which stores
    # prices/percentages as text, some containing thousand-separator dots or %% signs)
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell "D2" '29.054.80'
Set-TextCell "E2" '  -0.25%  '
Set-TextCell "D3" '1.817.75'
Set-TextCell "E3" '  -0.93%  '
Set-TextCell "D4" '0.9990'
Set-TextCell "E4" '  -0.31%  '
Set-TextCell "E5" '  -1.03%  '
Set-TextCell "D6" '0.6143'
Set-TextCell "E6" '  -2.20%  '
Set-TextCell "E7" '  -0.19%  '
Set-TextCell "D8" '0.07309'
Set-TextCell "E8" '  -2.51%  '
Set-TextCell "D9" '0.2875'
Set-TextCell "E9" '  -1.63%  '
Set-TextCell "D10" '22.81'
Set-TextCell "E10" '  -2.04%  '
Set-TextCell "E11" '  -0.46%  '
Set-TextCell "D12" '1.812.42'
Set-TextCell "E12" '  -1.22%  '
Set-TextCell "D13" '4.932'
Set-TextCell "E13" '  -1.56%  '
Set-TextCell "E14" '  -1.31%  '
Set-TextCell "D15" '81.43'
Set-TextCell "E15" '  -1.67%  '
Set-TextCell "D16" '0.000008997'
Set-TextCell "E16" '  -3.98%  '
Set-TextCell "E17" '  -2.65%  '
Set-TextCell "D18" '29.044.40'
Set-TextCell "E18" '  -0.30%  '
Set-TextCell "D19" '2.068.27'
Set-TextCell "E19" '  -0.65%  '
Set-TextCell "D20" '236.91'
Set-TextCell "E20" '  +6.12%  '
Set-TextCell "D21" '12.40'
Set-TextCell "E21" '  -1.55%  '
Set-TextCell "E22" '  -0.34%  '
Set-TextCell "D23" '7.098'
Set-TextCell "E23" '  +0.03%  '
Set-TextCell "D24" '1.000'
Set-TextCell "E24" '  -0.31%  '
Set-TextCell "D25" '157.64'
Set-TextCell "E25" '  -1.53%  '
Set-TextCell "D26" '0.1401'
Set-TextCell "E26" '  +0.67%  '
Set-TextCell "D27" '8.407'
Set-TextCell "E27" '  -1.09%  '
Set-TextCell "D28" '17.54'
Set-TextCell "E28" '  -2.01%  '
Set-TextCell "D29" '1.484'
Set-TextCell "E29" '  -0.96%  '
Set-TextCell "D30" '0.05559'
Set-TextCell "E30" '  -1.42%  '
Set-TextCell "E31" '  -0.16%  '
Set-TextCell "D32" '4.083'
Set-TextCell "E32" '  -1.76%  '
Set-TextCell "E33" '  -0.22%  '
Set-TextCell "D34" '0.7329'
Set-TextCell "E34" '  -1.21%  '
Set-TextCell "D35" '1.808'
Set-TextCell "E35" '  -1.72%  '
Set-TextCell "E36" '  -1.19%  '
Set-TextCell "D37" '2.615'
Set-TextCell "E37" '  -2.20%  '
Set-TextCell "D38" '2.824'
Set-TextCell "E38" '  +2.20%  '
Set-TextCell "D39" '1.206.11'
Set-TextCell "E39" '  -1.28%  '
Set-TextCell "E40" '  -1.44%  '
Set-TextCell "D41" '6.349'
Set-TextCell "E41" '  -3.01%  '
Set-TextCell "D42" '0.8908'
Set-TextCell "E42" '  -0.15%  '
Set-TextCell "D43" '1.001'
Set-TextCell "E43" '  -0.15%  '
Set-TextCell "D44" '100.82'
Set-TextCell "E44" '  -1.16%  '
Set-TextCell "D45" '1.971.62'
Set-TextCell "E45" '  -0.46%  '
Set-TextCell "D46" '64.26'
Set-TextCell "E46" '  -2.36%  '
Set-TextCell "D47" '0.5084'
Set-TextCell "E47" '  -0.28%  '
Set-TextCell "E48" '  -4.33%  '
Set-TextCell "D49" '0.3984'
Set-TextCell "E49" '  -2.26%  '
Set-TextCell "D50" '9.005'
Set-TextCell "E50" '  -0.02%  '
Set-TextCell "D51" '0.05746'
Set-TextCell "E51" '  -1.38%  '
